# Re-pull / push updated dSF (column F) data into the calendar sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new dSF (column F) value, as re-pulled from source data.
$newF = @{
    2  = -1
    3  = -1
    4  = -1
    5  = -1
    7  = 1
    8  = -1
    9  = 5
    10 = -1
    11 = 4
    12 = 7
    13 = -3
    14 = 8
    15 = 4
    17 = 1
    19 = 10
    20 = 4
    21 = -4
    24 = -4
    25 = 2
    26 = 2
    27 = 6
    28 = -2
    30 = -2
    31 = -2
    32 = -2
    33 = -1
    34 = 1
    35 = 11
    36 = -2
    37 = -1
}

foreach ($row in $newF.Keys) {
    $ws.Cells.Item($row, 6).Value = $newF[$row]
}
